# Insert a new weekly record row for "Feria Lagunitas de Puerto Montt - Apio"
# at sheet row 455, pushing the existing rows 455..538 down to 456..539.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 455 (Excel shifts row 455..538 down to 456..539,
# copying formatting - e.g. the date number format on column D - from the row above).
$ws.Rows.Item(455).Insert()

# Populate the newly inserted row with the new market record.
$ws.Cells.Item(455, 1).Value = 4
$ws.Cells.Item(455, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(455, 3).Value = "Los Lagos"
$ws.Cells.Item(455, 4).Value = 45209
$ws.Cells.Item(455, 5).Value = 10
$ws.Cells.Item(455, 6).Value = 100112017
$ws.Cells.Item(455, 7).Value = "Apio"
$ws.Cells.Item(455, 8).Value = "Americana (o)"
$ws.Cells.Item(455, 9).Value = "Primera"
$ws.Cells.Item(455, 10).Value = 45
$ws.Cells.Item(455, 11).Value = 11000
$ws.Cells.Item(455, 12).Value = 11000
$ws.Cells.Item(455, 13).Value = 11000
$ws.Cells.Item(455, 14).Value = "$/docena de matas"
$ws.Cells.Item(455, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(455, 16).Value = 1833
$ws.Cells.Item(455, 17).Value = 6
$ws.Cells.Item(455, 18).Value = "Hortaliza"
